$d = $word.ActiveDocument

# wdReplaceAll = 2, wdFindContinue = 1

# 1) Fix the typo "Test Driven Develop (TDD)" -> "Test Driven Development (TDD)"
$d.Content.Find.Execute(
    "Test Driven Develop (TDD) is the act", $true, $false, $false, $false, $false,
    $true, 1, $false, "Test Driven Development (TDD) is the act", 2)

# 2) "functionality is working, these tests" -> "functionality is working; These tests"
$d.Content.Find.Execute(
    "functionality is working, these tests", $true, $false, $false, $false, $false,
    $true, 1, $false, "functionality is working; These tests", 2)

# 3) Mention Selenium: "unit tests, but can also" -> "unit tests using libraries such as Selenium, but can also"
$d.Content.Find.Execute(
    "like unit tests, but can also", $true, $false, $false, $false, $false,
    $true, 1, $false, "like unit tests using libraries such as Selenium, but can also", 2)

# 4) "a big part of the test methodology" -> "a big part of the TDD methodology"
$d.Content.Find.Execute(
    "a big part of the test methodology", $true, $false, $false, $false, $false,
    $true, 1, $false, "a big part of the TDD methodology", 2)

# 5) Merge the split "Tes" + "te" + "r: " runs in the test-plan table into one run "Tester: "
$d.Content.Find.Execute(
    "Tester: ", $true, $false, $false, $false, $false,
    $true, 1, $false, "Tester: ", 2)

# 6) Move the _GoBack bookmark from the start of the "Test Driven Development"
#    heading to the very end of the "Alongside TDD..." paragraph (after the
#    last run, before the paragraph mark).
$endRng = $d.Content
$endRng.Find.Execute("Alongside TDD")
$paraRng = $endRng.Paragraphs(1).Range
$paraRng.MoveEnd(1, -1)
$paraRng.Collapse(0)
$paraRng.InsertAfter("ZZ_BOOKMARK_MARKER_ZZ")

$markerRng = $d.Content
$markerRng.Find.Execute("ZZ_BOOKMARK_MARKER_ZZ")
$markerRng.Collapse(1)
$d.Bookmarks.Add("_GoBack", $markerRng)

$cleanupRng = $d.Content
$cleanupRng.Find.Execute("ZZ_BOOKMARK_MARKER_ZZ")
$cleanupRng.Text = ""
